# RPA datasets push 2024-08-09
# Update confirmed-offering-price (D column) and offering-amount-in-millions
# (E column) figures for a few IPO rows now that the results became known.
#
# These figures are stored as text (shared strings) even though they look
# like plain numbers, so each cell is briefly switched to a text number
# format ("@") before the assignment, then the style is reset back to
# "Normal" afterwards so no left-over quote-prefix / custom number format
# is attached to the cell (matching the original workbook's styling).

function Set-TextValue {
    param($worksheet, $Address, $Text)
    $rng = $worksheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 대신밸런스스팩18호 (row 7): 확정공모가 "-" -> "2000"
Set-TextValue $ws "D7" "2000"

# 티디에스팜 (row 10): 확정공모가 "-" -> "13000", 공모금액(백만) "9500" -> "13000"
Set-TextValue $ws "D10" "13000"
Set-TextValue $ws "E10" "13000"

# 케이쓰리아이 (row 11): 공모금액(백만) "22351" -> "21700"
Set-TextValue $ws "E11" "21700"
